$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Duplicate the class-subset table (header row2, 19 data rows, blank, count row23) ---
# --- down to rows 25-46, bringing along values and cell formatting (fills, borders). ---
$ws.Range("B2:K23").Copy($ws.Range("B25"))

# --- 2. The copy above duplicated the old boolean flags verbatim. Every cell that was TRUE ---
# --- now needs to show the running "how many classes match so far" count instead; cells ---
# --- that were FALSE keep showing the red FALSE flag, unchanged. ---

function Set-ColumnValues($sheet, $rangeAddr, $values) {
    $rng = $sheet.Range($rangeAddr)
    $arr = New-Object 'object[,]' $values.Count,1
    for ($i = 0; $i -lt $values.Count; $i++) { $arr[$i,0] = $values[$i] }
    $rng.Value = $arr
}

# Column D
Set-ColumnValues $ws "D26:D44" @(0,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18)

# Column E
Set-ColumnValues $ws "E26:E34" @(0,1,2,3,4,5,6,7,8)
Set-ColumnValues $ws "E36:E41" @(9,10,11,12,13,14)
Set-ColumnValues $ws "E43:E44" @(15,16)

# Column F
Set-ColumnValues $ws "F26:F39" @(0,1,2,3,4,5,6,7,8,9,10,11,12,13)
Set-ColumnValues $ws "F43:F44" @(14,15)

# Column G
Set-ColumnValues $ws "G26:G34" @(0,1,2,3,4,5,6,7,8)
Set-ColumnValues $ws "G36:G39" @(9,10,11,12)
$ws.Range("G41").Value = 13
Set-ColumnValues $ws "G43:G44" @(14,15)

# Column H
Set-ColumnValues $ws "H26:H28" @(0,1,2)
Set-ColumnValues $ws "H32:H34" @(3,4,5)
Set-ColumnValues $ws "H36:H39" @(6,7,8,9)
$ws.Range("H41").Value = 10
Set-ColumnValues $ws "H43:H44" @(11,12)

# --- 3. J (idda+synthia) and K (idda+cc) combine two/three of the other flags with AND(). ---
# --- Re-derive: where the combination is still FALSE, keep it a live formula; where it is ---
# --- now TRUE, show the running count like the other columns. ---
$ws.Range("K29").Formula = "=AND(F29,H29)"
$ws.Range("K30").Formula = "=AND(F30,H30)"
$ws.Range("K31").Formula = "=AND(F31,H31)"
$ws.Range("J35").Formula = "=AND(F35:G35)"
$ws.Range("K35").Formula = "=AND(F35,H35)"
$ws.Range("J40").Formula = "=AND(F40:G40)"
$ws.Range("K40").Formula = "=AND(F40,H40)"
$ws.Range("J41").Formula = "=AND(F41:G41)"
$ws.Range("K41").Formula = "=AND(F41,H41)"
$ws.Range("J42").Formula = "=AND(F42:G42)"
$ws.Range("K42").Formula = "=AND(F42,H42)"

# Column J (remaining TRUE cells)
Set-ColumnValues $ws "J26:J34" @(0,1,2,3,4,5,6,7,8)
Set-ColumnValues $ws "J36:J39" @(9,10,11,12)
Set-ColumnValues $ws "J43:J44" @(13,14)

# Column K (remaining TRUE cells)
Set-ColumnValues $ws "K26:K28" @(0,1,2)
Set-ColumnValues $ws "K32:K34" @(3,4,5)
Set-ColumnValues $ws "K36:K39" @(6,7,8,9)
Set-ColumnValues $ws "K43:K44" @(10,11)

# --- 4. The original colour-scale conditional formatting only made sense for the single ---
# --- boolean table; now that results are reported for every subset it is dropped. ---
$ws.Range("D3:H21").FormatConditions.Delete()

# --- 5. Leave the selection where the author finished editing, at the bottom-right of the ---
# --- newly added table. ---
$ws.Range("K45").Select()
